$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph whose text contains a given substring.
# ---------------------------------------------------------------------------
function Get-ParagraphContaining {
    param($doc, [string]$needle)
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. "As a fighter, I …" -> "As a fighter, I have a small chance to dodge an
#    incoming attack", and strike the whole bullet through (it is being cut
#    from the backlog).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("As a fighter, I …", $true, $false, $false, $false, `
                         $false, $true, 1, $false, `
                         "As a fighter, I have a small chance to dodge an incoming attack", `
                         2) | Out-Null

$fighterPara = Get-ParagraphContaining $d "As a fighter, I have a small chance to dodge an incoming attack"
$fighterPara.Range.Font.StrikeThrough = $true

# ---------------------------------------------------------------------------
# 2. Word always keeps a single "_GoBack" bookmark marking the location of
#    the most recent edit. Since the last edit of this revision happened
#    right after "no points", move the bookmark there (removing it from
#    wherever it previously sat, at the very end of the document).
# ---------------------------------------------------------------------------
$pointsPara = Get-ParagraphContaining $d "no points"

# Position right after the last real character of the paragraph (i.e. before
# its paragraph mark): Range.End sits just past the paragraph mark, so back
# up one.
$insertPos = $pointsPara.Range.End - 1

$ip = $pointsPara.Range
$ip.Start = $insertPos
$ip.End = $insertPos

# Use a throwaway character so the bookmark range is non-collapsed while we
# create it (collapsed ranges exactly on a paragraph boundary are resolved
# unreliably), then delete the character again, leaving the bookmark
# correctly anchored immediately after "no points".
$ip.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $ip)

$cleanup = $pointsPara.Range
$cleanup.Start = $insertPos
$cleanup.End = $insertPos + 1
$cleanup.Delete()
